# emprestimo.xlsx - update loan log: fix technician assignments on a few
# existing rows, and register the new loan entries that were handed out
# since the last save.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# ---------------------------------------------------------------------------
# Fixes to existing rows (row 5, 7, 8, 9): correct the "Pointer" name typo
# and re-assign a few loans from Kauan to Victor, plus flip their delivery
# status now that the gear has actually been handed over.
# ---------------------------------------------------------------------------
$ws.Range("A5").Value = "Pointer"
$ws.Range("F5").Value = "Entregue "

$ws.Range("D7").Value = "Victor"

$ws.Range("D8").Value = "Victor"
$ws.Range("F8").Value = "Entregue"

$ws.Range("D9").Value = "Victor"
$ws.Range("F9").Value = "Entregue"

# ---------------------------------------------------------------------------
# New loan rows (11-17)
# ---------------------------------------------------------------------------
function Set-Row($r, $equip, $sala, $y, $m, $d, $tecnico, $responsavel, $status) {
    $ws.Cells.Item($r, 1).Value = $equip
    $ws.Cells.Item($r, 2).Value = $sala
    $ws.Cells.Item($r, 3).Value = (Get-Date -Year $y -Month $m -Day $d -Hour 0 -Minute 0 -Second 0).Date
    $ws.Cells.Item($r, 4).Value = $tecnico
    $ws.Cells.Item($r, 5).Value = $responsavel
    $ws.Cells.Item($r, 6).Value = $status
}

Set-Row 11 "Notebook Vaio" "Extra Aula" 2024 9 30 "Victor" "Aluna" "Entregue"
Set-Row 12 "Microfone" "219B" 2024 10 1 "Victor" "Prof" "Entregue"
Set-Row 13 "Microfone" "322B" 2024 10 1 "Victor" "Prof" "Entregue"
Set-Row 14 "Pointer" "316D" 2024 10 3 "Kauan" "Thaise Gerber" "Entregue"
Set-Row 15 "Microfone" "219B" 2024 10 8 "Victor" "Prof" "Não entregue"
Set-Row 16 "Notebook Vaio" "211H" 2024 10 16 "Victor" "Nathalia (48 999467019)" "Entregue "
Set-Row 17 "Caixa de Som" "318D" 2024 10 16 "Victor" "Jurema (48 99962-9662)" "Não entregue"

# Row 15's status cell gets a distinct look: underlined, to flag that this
# one hasn't come back yet. (Border/centering are already inherited from
# the row's existing style.)
$ws.Range("F15").Font.Underline = 2

# ---------------------------------------------------------------------------
# Cursor was left on F18 at last save.
# ---------------------------------------------------------------------------
$ws.Range("F18").Select()
